# PolicyEngine budget workbook update
# Rewrites "FY2026 Budget" sheet into "2025 Budget" with actual
# spending data through October 2025, per commit message / diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Rename the sheet
# ---------------------------------------------------------------
$ws.Name = "2025 Budget"

# ---------------------------------------------------------------
# 2. Park copies of the existing styled cells (format only) in a
#    scratch area far to the right so we can reapply the exact
#    same cell styles (fonts/fills/number formats) after we wipe
#    and rebuild the main A1:C31 area. This avoids creating
#    duplicate style entries in the workbook.
# ---------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # style1: bold 16  (title)

$ws.Range("A2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # style2: bold 12  (subtitle)

$ws.Range("A3").Copy()
$ws.Range("Z3").PasteSpecial(-4122)   # style3: italic 10 (data source line)

$ws.Range("A5").Copy()
$ws.Range("Z4").PasteSpecial(-4122)   # style4: bold (section label / total label)

$ws.Range("B8").Copy()
$ws.Range("Z5").PasteSpecial(-4122)   # style5: currency $#,##0

$ws.Range("A7").Copy()
$ws.Range("Z6").PasteSpecial(-4122)   # style6: bold12 + fill (section header)

$ws.Range("B12").Copy()
$ws.Range("Z7").PasteSpecial(-4122)   # style7: currency bold (grand total)

# ---------------------------------------------------------------
# 3. Clear the old body content (formats + values), leaving our
#    parked scratch cells (column Z) untouched.
# ---------------------------------------------------------------
$ws.Range("A1:C26").Clear()

# ---------------------------------------------------------------
# 4. Column widths (40 / 18 / 12 chars as stored in the XML). The
#    COM ColumnWidth property is offset from the raw XML width by
#    ~0.8333 characters (standard Calibri 11 padding), so we back
#    that out to land exactly on the target stored widths.
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 39.16666666666667
$ws.Columns.Item(2).ColumnWidth = 17.16666666666667
$ws.Columns.Item(3).ColumnWidth = 11.16666666666667

# ---------------------------------------------------------------
# 5. Header block
# ---------------------------------------------------------------
$ws.Range("A1").Value = "PolicyEngine"
$ws.Range("Z1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("A2").Value = "Calendar Year 2025 Budget (Actual Spending Through October)"
$ws.Range("Z2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A3").Value = "Data source: opencollective.com/policyengine (348 transactions analyzed)"
$ws.Range("Z3").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 6. Summary block (rows 5-7)
# ---------------------------------------------------------------
$ws.Range("A5").Value = "Year-to-Date Actuals (Jan-Oct 2025)"
$ws.Range("Z4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").Value = 470611
$ws.Range("Z5").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("A6").Value = "Annualized Projection (12 months)"
$ws.Range("Z4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").Value = 564734
$ws.Range("Z5").Copy()
$ws.Range("B6").PasteSpecial(-4122)

$ws.Range("A7").Value = "Current Balance (Oct 2025)"
$ws.Range("B7").Value = 545558
$ws.Range("Z5").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 7. Expense breakdown block (rows 9-23)
# ---------------------------------------------------------------
$ws.Range("A9").Value = "EXPENSE BREAKDOWN (Annualized)"
$ws.Range("Z6").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$expenseRows = @(
    @(10, "Fellowships & Internships",    103904, "18.4%"),
    @(11, "PolicyEngine UK Operations",   100965, "17.9%"),
    @(12, "Contractor Services",           94170, "16.7%"),
    @(13, "Cloud Infrastructure",          68868, "12.2%"),
    @(14, "Website Development",           36000, "6.4%"),
    @(15, "Travel & Conferences",          28751, "5.1%"),
    @(16, "Payroll Expenses",              22021, "3.9%"),
    @(17, "Legal Services",                13878, "2.5%"),
    @(18, "AI/Software Subscriptions",     10145, "1.8%"),
    @(19, "Equipment & Materials",          7994, "1.4%"),
    @(20, "Office Space",                   6672, "1.2%"),
    @(21, "Events & Supplies",              1683, "0.3%"),
    @(22, "Other Operating",               69683, "12.3%")
)

# Force the percent column to stay as literal text (not an
# auto-converted percentage number) from the very first write,
# so no stray percentage-number style is ever created.
$ws.Range("C10:C22").NumberFormat = "@"

foreach ($row in $expenseRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("Z5").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("C$r").Value = $row[3]
}

# Reset the percent column back to the default (unstyled) cell
# style, matching the source workbook, while keeping the literal
# text values already entered.
$ws.Range("C10:C22").Style = "Normal"

$ws.Range("A23").Value = "TOTAL EXPENSES (Annualized)"
$ws.Range("Z4").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("B23").Value = 564734
$ws.Range("Z7").Copy()
$ws.Range("B23").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 8. Revenue block (rows 25-26)
# ---------------------------------------------------------------
$ws.Range("A25").Value = "REVENUE (YTD Jan-Oct 2025)"
$ws.Range("Z6").Copy()
$ws.Range("A25").PasteSpecial(-4122)

$ws.Range("A26").Value = "Foundation Grants & Contributions"
$ws.Range("B26").Value = 643910
$ws.Range("Z5").Copy()
$ws.Range("B26").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 9. Note block (rows 28-31), new styles:
#    - A28 "Note:" -> italic, default size/font (new font, no size)
#    - A29 long note -> default font, wrap text + vertical-top
#      alignment, merged across A29:C31
# ---------------------------------------------------------------
$ws.Range("A28").Value = "Note:"
$ws.Range("A28").Font.Italic = $true

$ws.Range("A29").Value = 'Open Collective shows estimated annual budget of $1,100,315, but actual spending through October 2025 totals $470,611 (10 months), projecting to ~$565,000 annualized. This budget reflects actual spending patterns. The difference may reflect planned but not yet executed activities or conservative budgeting.'

# Merge first (while A29 still has the default style) so the
# wrap/vertical-top formatting applied next stays local to A29
# instead of spreading across the whole merged range.
$ws.Range("A29:C31").Merge()

$ws.Range("A29").WrapText = $true
$ws.Range("A29").VerticalAlignment = -4160

# Materialize empty row stubs for rows 30 and 31 (touched by the
# merge range but otherwise empty), matching the source workbook.
$ws.Rows.Item(30).OutlineLevel = 0
$ws.Rows.Item(31).OutlineLevel = 0

# ---------------------------------------------------------------
# 10. Clean up scratch/parking cells used for format copying.
# ---------------------------------------------------------------
$ws.Range("Z1:Z10").Clear()

# ---------------------------------------------------------------
# 11. Leave active selection at A1, like the original workbook.
# ---------------------------------------------------------------
$ws.Range("A1").Select()
